# Adds RF035, RF036 and RF037 requirement paragraphs right after RF034
# ("O sistema deve permitir que organizador do evento anexe arquivos
# durante a criação do evento.") and before the "Requisitos Não
# Funcionais" heading.

$d = $word.ActiveDocument

# Locate the last paragraph of RF034 ("O sistema deve permitir ...").
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "O sistema deve permitir que organizador do evento anexe arquivos durante a criação do evento.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the RF034 paragraph to anchor the insertion."
}

$anchorParagraph = $searchRange.Paragraphs(1)

# Collapse to just before that paragraph's end-of-paragraph mark so the
# new content is inserted inside the existing BodyText paragraph range,
# which keeps the new paragraphs' style as BodyText without having to
# stamp an explicit paragraph-style change (and its rsid side effect).
$insertPos = $anchorParagraph.Range.End - 1
$cursor = $d.Range($insertPos, $insertPos)

$openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">RF035: Verificação de email da conta de usuário</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve">–</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Ademir, Gustavo, Matheus</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:t xml:space="preserve">Na execução do cadastro o sistema deve enviar um email de verificação com um link para confirmar que o email pertence a pessoa. Ao seguir o link a conta do usuário será ativada.</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">RF036: Sugestões de provedores de email</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve">–</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Ademir, Gustavo</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:t xml:space="preserve">Todos os campos de email do sistema devem oferecer sugestões de provedores de email. Os provedores devem ser sugeridos quando o usuário digitar @.</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">RF037: Impressão da lista de convidados</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve">–</w:t></w:r>' + `
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
      '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Gustavo</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + `
      '<w:r><w:t xml:space="preserve">Na página da lista de convidados deve ter uma opção para realizar a impressão da lista com os nomes e número de acompanhantes confirmados de todos os convidados confirmados.</w:t></w:r>' + `
    '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData>' + `
  '</pkg:part>' + `
  '</pkg:package>'

$cursor.InsertXML($openXml) | Out-Null

Write-Output "Inserted RF035, RF036 and RF037."
